# Apply "roll forward one fiscal year" update + new publish dates, per commit:
# "update database and change read_price algorithm"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row 8: fiscal period labels (columns D:H) ---
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# --- Header row 9: publish dates (columns D:H) ---
$ws.Range("D9").Value = "1399-12-18 (3)"
$ws.Range("E9").Value = "1400-04-14 (8)"
$ws.Range("F9").Value = "1401-04-12 (11)"
$ws.Range("G9").Value = "1402-02-25 (8)"
# H9 = "1402-02-25" looks exactly like a date (yyyy-mm-dd), so a plain .Value
# assignment gets auto-converted to a date serial and also forces a brand new
# number-format style onto the cell. Enter it as a text formula instead (so
# it is stored verbatim, not date-parsed) and then convert that formula back
# down to a plain value in-place; this keeps the original "General" style
# (s="9") intact and avoids creating any stray style entries.
$ws.Range("H9").Formula = "=""1402-02-25"""
$ws.Range("H9").Copy() | Out-Null
$ws.Range("H9").PasteSpecial(-4163) | Out-Null

# --- Row 11: فروش (Sales) ---
$ws.Range("D11").Value = 3003011
$ws.Range("E11").Value = 6567810
$ws.Range("F11").Value = 12626079
$ws.Range("G11").Value = 17210109
$ws.Range("H11").Value = 27260727

# --- Row 12: بهای تمام شده کالای فروش رفته (Cost of goods sold) ---
$ws.Range("D12").Value = -1668107
$ws.Range("E12").Value = -3272475
$ws.Range("F12").Value = -6467983
$ws.Range("G12").Value = -9208790
$ws.Range("H12").Value = -15164059

# --- Row 13: سود (زیان) ناخالص (Gross profit) ---
$ws.Range("D13").Value = 1334904
$ws.Range("E13").Value = 3295335
$ws.Range("F13").Value = 6158096
$ws.Range("G13").Value = 8001319
$ws.Range("H13").Value = 12096668

# --- Row 14: هزینه های عمومی, اداری و تشکیلاتی (General & admin expenses) ---
$ws.Range("D14").Value = -800082
$ws.Range("E14").Value = -1096680
$ws.Range("F14").Value = -1790320
$ws.Range("G14").Value = -3084627
$ws.Range("H14").Value = -4765322

# --- Row 15: هزینه کاهش ارزش دریافتنی‌ها (D15 stays "-", E:H shift) ---
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0

# --- Row 16: خالص سایر درامدها (هزینه ها) ی عملیاتی ---
$ws.Range("D16").Value = 352
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 0

# --- Row 17: سود (زیان) عملیاتی (Operating profit) ---
$ws.Range("D17").Value = 535174
$ws.Range("E17").Value = 2198655
$ws.Range("F17").Value = 4367776
$ws.Range("G17").Value = 4916692
$ws.Range("H17").Value = 7331346

# --- Row 18: هزینه های مالی (Financial expenses) ---
$ws.Range("D18").Value = -322094
$ws.Range("E18").Value = -413983
$ws.Range("F18").Value = -516842
$ws.Range("G18").Value = -674205
$ws.Range("H18").Value = -1038607

# --- Row 19: خالص سایر درامدها و هزینه های غیرعملیاتی ---
$ws.Range("D19").Value = -99548
$ws.Range("E19").Value = -36080
$ws.Range("F19").Value = -49793
$ws.Range("G19").Value = 97616
$ws.Range("H19").Value = 101323

# --- Row 20: سود (زیان) خالص عملیات در حال تداوم قبل از مالیات ---
$ws.Range("D20").Value = 113532
$ws.Range("E20").Value = 1748592
$ws.Range("F20").Value = 3801141
$ws.Range("G20").Value = 4340103
$ws.Range("H20").Value = 6394062

# --- Row 21: مالیات (Tax) ---
$ws.Range("D21").Value = -24472
$ws.Range("E21").Value = -247384
$ws.Range("F21").Value = -192068
$ws.Range("G21").Value = -217880
$ws.Range("H21").Value = -256002

# --- Row 22: سود (زیان) خالص عملیات در حال تداوم ---
$ws.Range("D22").Value = 89060
$ws.Range("E22").Value = 1501208
$ws.Range("F22").Value = 3609073
$ws.Range("G22").Value = 4122223
$ws.Range("H22").Value = 6138060

# --- Row 23: سود (زیان) عملیات متوقف شده پس از اثر مالیاتی (stays all zero) ---
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0

# --- Row 24: سود (زیان) خالص ---
$ws.Range("D24").Value = 89060
$ws.Range("E24").Value = 1501208
$ws.Range("F24").Value = 3609073
$ws.Range("G24").Value = 4122223
$ws.Range("H24").Value = 6138060

# --- Row 25: سود هر سهم پس از کسر مالیات ---
$ws.Range("D25").Value = 178
$ws.Range("E25").Value = 3002
$ws.Range("F25").Value = 2406
$ws.Range("G25").Value = 2748
$ws.Range("H25").Value = 1535

# --- Row 26: سرمایه (Capital) ---
$ws.Range("D26").Value = 500000
$ws.Range("E26").Value = 500000
$ws.Range("F26").Value = 1500000
$ws.Range("G26").Value = 1500000
$ws.Range("H26").Value = 4000000

# --- Row 27: سود هر سهم بر اساس آخرین سرمایه ---
$ws.Range("D27").Value = 22
$ws.Range("E27").Value = 375
$ws.Range("F27").Value = 902
$ws.Range("G27").Value = 1031
$ws.Range("H27").Value = 1535

$wb.Save()
